# Append a new log row (row 83) to each of the four worksheets.
# The logging system used to append one row per day; its configuration is
# being removed, but the final row it produced (for the latest day) still
# needs to be captured here. Each new row duplicates the previous last row
# (row 82) verbatim except for the timestamp in column A, which advances to
# the next day's snapshot time.

$wb = $excel.ActiveWorkbook

$newTime = 45869.49581018519

foreach ($idx in 1, 2, 3, 4) {
    $ws = $wb.Worksheets.Item($idx)

    # Duplicate the last existing row (82) into the new row (83), carrying
    # over formatting/styles (e.g. the date style on column A) exactly.
    $ws.Range("A82:I82").Copy($ws.Range("A83:I83"))

    # The new row records the next day's log snapshot.
    $ws.Cells.Item(83, 1).Value = $newTime
}
